# Auto-generated edit script: fills in confusion-matrix experiment values
# that were previously blank (t-SNE MI_2000 LR_700 PE_16 results workbook).
# Use manual calculation so the already-cached #DIV/0! formula results
# (F and M columns) are left untouched, matching the target diff exactly.
$excel.Calculation = -4135  # xlCalculationManual

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FS")
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 47.368421052631575
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 7.017543859649122
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.013888888888888888
$ws.Range("K25").Value = 0.9166666666666666
$ws.Range("L25").Value = 0.06944444444444445
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.8333333333333334
$ws.Range("E26").Value = 0.0
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.9074074074074074

$ws = $wb.Worksheets.Item("IF")
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 42.10526315789473
$ws.Range("J24").Value = 0.8444444444444444
$ws.Range("K24").Value = 0.15555555555555556
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 9.941520467836257
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.625
$ws.Range("E25").Value = 0.375
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9583333333333334
$ws.Range("L25").Value = 0.041666666666666664
$ws.Range("C26").Value = 0.3333333333333333
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.3333333333333333
$ws.Range("J26").Value = 0.05555555555555555
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.8703703703703703

$ws = $wb.Worksheets.Item("IA")
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.022222222222222223
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 5.847953216374268
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.875
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.3333333333333333
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.16666666666666666
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.1111111111111111
$ws.Range("L26").Value = 0.8703703703703703

$ws = $wb.Worksheets.Item("FS-IF")
$ws.Range("C24").Value = 0.4
$ws.Range("D24").Value = 0.4
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 31.57894736842105
$ws.Range("J24").Value = 1.0
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 4.678362573099415
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.1111111111111111
$ws.Range("L26").Value = 0.8888888888888888

$ws = $wb.Worksheets.Item("FS-IA")
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 31.57894736842105
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 2.923976608187134
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.5
$ws.Range("E25").Value = 0.5
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9861111111111112
$ws.Range("L25").Value = 0.013888888888888888
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.05555555555555555
$ws.Range("L26").Value = 0.9444444444444444

$ws = $wb.Worksheets.Item("IF-IA")
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 21.052631578947366
$ws.Range("J18").Value = 1.0
$ws.Range("K18").Value = 0.0
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 0.0
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.875
$ws.Range("E19").Value = 0.0
$ws.Range("J19").Value = 0.0
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.0
$ws.Range("L20").Value = 1.0
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 15.789473684210526
$ws.Range("J24").Value = 0.9555555555555556
$ws.Range("K24").Value = 0.022222222222222223
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 8.187134502923977
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.875
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9166666666666666
$ws.Range("L25").Value = 0.08333333333333333
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.8888888888888888

$ws = $wb.Worksheets.Item("FS-IF-IA")
$ws.Range("C18").Value = 0.6
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.2
$ws.Range("G18").Value = 26.31578947368421
$ws.Range("J18").Value = 1.0
$ws.Range("K18").Value = 0.0
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 1.1695906432748537
$ws.Range("C19").Value = 0.0
$ws.Range("D19").Value = 0.875
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.013888888888888888
$ws.Range("K19").Value = 0.9861111111111112
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.018518518518518517
$ws.Range("L20").Value = 0.9814814814814815
$ws.Range("C24").Value = 0.4
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.6
$ws.Range("G24").Value = 31.57894736842105
$ws.Range("J24").Value = 1.0
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 6.432748538011696
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.875
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.013888888888888888
$ws.Range("K25").Value = 0.9305555555555556
$ws.Range("L25").Value = 0.05555555555555555
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.09259259259259259
$ws.Range("L26").Value = 0.8888888888888888

Write-Output "done"